$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 9,16
$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 12.39269
$arr[0,3] = 37.17807
$arr[0,4] = 0.6136856461363306
$arr[0,5] = 0.6136856461363306
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 5.188906999999999
$arr[0,9] = 15.566721
$arr[0,10] = 0.02571200377994867
$arr[0,11] = 0.02571200377994868
$arr[0,12] = 64.30451588983
$arr[0,13] = 578.74064300847
$arr[0,14] = 0.01577908765315757
$arr[0,15] = 0.01577908765315758
$arr[1,0] = 3
$arr[1,1] = 1
$arr[1,2] = 12.39269
$arr[1,3] = 37.17807
$arr[1,4] = 0.6136856461363306
$arr[1,5] = 0.6136856461363306
$arr[1,6] = 3
$arr[1,7] = 1
$arr[1,8] = 67.633555
$arr[1,9] = 202.900665
$arr[1,10] = 0.3351369029761694
$arr[1,11] = 0.3351369029761695
$arr[1,12] = 838.16168071295
$arr[1,13] = 7543.455126416549
$arr[1,14] = 0.2056687068470593
$arr[1,15] = 0.2056687068470593
$arr[2,0] = 3
$arr[2,1] = 1
$arr[2,2] = 12.39269
$arr[2,3] = 37.17807
$arr[2,4] = 0.6136856461363306
$arr[2,5] = 0.6136856461363306
$arr[2,6] = 3
$arr[2,7] = 1
$arr[2,8] = 128.9862746666667
$arr[2,9] = 386.958824
$arr[2,10] = 0.6391510932438819
$arr[2,11] = 0.6391510932438819
$arr[2,12] = 1598.486916198853
$arr[2,13] = 14386.38224578968
$arr[2,14] = 0.3922378516361137
$arr[2,15] = 0.3922378516361137
$arr[3,0] = 3
$arr[3,1] = 1
$arr[3,2] = 3.644292666666666
$arr[3,3] = 10.932878
$arr[3,4] = 0.1804652662055796
$arr[3,5] = 0.1804652662055796
$arr[3,6] = 3
$arr[3,7] = 1
$arr[3,8] = 5.188906999999999
$arr[3,9] = 15.566721
$arr[3,10] = 0.02571200377994867
$arr[3,11] = 0.02571200377994868
$arr[3,12] = 18.90989572811533
$arr[3,13] = 170.189061553038
$arr[3,14] = 0.004640123606827306
$arr[3,15] = 0.004640123606827307
$arr[4,0] = 3
$arr[4,1] = 1
$arr[4,2] = 3.644292666666666
$arr[4,3] = 10.932878
$arr[4,4] = 0.1804652662055796
$arr[4,5] = 0.1804652662055796
$arr[4,6] = 3
$arr[4,7] = 1
$arr[4,8] = 67.633555
$arr[4,9] = 202.900665
$arr[4,10] = 0.3351369029761694
$arr[4,11] = 0.3351369029761695
$arr[4,12] = 246.4764685070967
$arr[4,13] = 2218.28821656387
$arr[4,14] = 0.06048057041090792
$arr[4,15] = 0.06048057041090794
$arr[5,0] = 3
$arr[5,1] = 1
$arr[5,2] = 3.644292666666666
$arr[5,3] = 10.932878
$arr[5,4] = 0.1804652662055796
$arr[5,5] = 0.1804652662055796
$arr[5,6] = 3
$arr[5,7] = 1
$arr[5,8] = 128.9862746666667
$arr[5,9] = 386.958824
$arr[5,10] = 0.6391510932438819
$arr[5,11] = 0.6391510932438819
$arr[5,12] = 470.0637348683858
$arr[5,13] = 4230.573613815472
$arr[5,14] = 0.1153445721878444
$arr[5,15] = 0.1153445721878444
$arr[6,0] = 3
$arr[6,1] = 1
$arr[6,2] = 4.156890333333333
$arr[6,3] = 12.470671
$arr[6,4] = 0.2058490876580898
$arr[6,5] = 0.2058490876580898
$arr[6,6] = 3
$arr[6,7] = 1
$arr[6,8] = 5.188906999999999
$arr[6,9] = 15.566721
$arr[6,10] = 0.02571200377994867
$arr[6,11] = 0.02571200377994868
$arr[6,12] = 21.56971734886566
$arr[6,13] = 194.127456139791
$arr[6,14] = 0.005292792519963791
$arr[6,15] = 0.005292792519963793
$arr[7,0] = 3
$arr[7,1] = 1
$arr[7,2] = 4.156890333333333
$arr[7,3] = 12.470671
$arr[7,4] = 0.2058490876580898
$arr[7,5] = 0.2058490876580898
$arr[7,6] = 3
$arr[7,7] = 1
$arr[7,8] = 67.633555
$arr[7,9] = 202.900665
$arr[7,10] = 0.3351369029761694
$arr[7,11] = 0.3351369029761695
$arr[7,12] = 281.1452709884683
$arr[7,13] = 2530.307438896215
$arr[7,14] = 0.06898762571820224
$arr[7,15] = 0.06898762571820227
$arr[8,0] = 3
$arr[8,1] = 1
$arr[8,2] = 4.156890333333333
$arr[8,3] = 12.470671
$arr[8,4] = 0.2058490876580898
$arr[8,5] = 0.2058490876580898
$arr[8,6] = 3
$arr[8,7] = 1
$arr[8,8] = 128.9862746666667
$arr[8,9] = 386.958824
$arr[8,10] = 0.6391510932438819
$arr[8,11] = 0.6391510932438819
$arr[8,12] = 536.1817982945449
$arr[8,13] = 4825.636184650903
$arr[8,14] = 0.1315686694199238
$arr[8,15] = 0.1315686694199238

$ws.Range("E2:T10").Value = $arr
